# ajuste: corrigindo as categorias
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total" column header (column X, after the last existing "Idade ignorada" column W)
$ws.Range("X1").Value = "Total"

# New per-row totals in column X for the existing disease-category rows (2-6)
$ws.Range("X2").Value = 2523
$ws.Range("X3").Value = 387
$ws.Range("X4").Value = 1059
$ws.Range("X5").Value = 380
$ws.Range("X6").Value = 1574

# New row 7: "Outros" category with its age-bracket breakdown and total
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 133
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 22
$ws.Range("F7").Value = 62
$ws.Range("G7").Value = 86
$ws.Range("H7").Value = 76
$ws.Range("I7").Value = 87
$ws.Range("J7").Value = 106
$ws.Range("K7").Value = 120
$ws.Range("L7").Value = 124
$ws.Range("M7").Value = 162
$ws.Range("N7").Value = 200
$ws.Range("O7").Value = 220
$ws.Range("P7").Value = 270
$ws.Range("Q7").Value = 289
$ws.Range("R7").Value = 348
$ws.Range("S7").Value = 378
$ws.Range("T7").Value = 286
$ws.Range("U7").Value = 122
$ws.Range("V7").Value = 33
$ws.Range("W7").Value = 7
$ws.Range("X7").Value = 3147

# New row 8: "Total" row summing all categories, with its age-bracket breakdown and overall total
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 151
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 14
$ws.Range("E8").Value = 33
$ws.Range("F8").Value = 74
$ws.Range("G8").Value = 104
$ws.Range("H8").Value = 119
$ws.Range("I8").Value = 153
$ws.Range("J8").Value = 210
$ws.Range("K8").Value = 284
$ws.Range("L8").Value = 347
$ws.Range("M8").Value = 531
$ws.Range("N8").Value = 712
$ws.Range("O8").Value = 840
$ws.Range("P8").Value = 963
$ws.Range("Q8").Value = 1017
$ws.Range("R8").Value = 1104
$ws.Range("S8").Value = 1117
$ws.Range("T8").Value = 861
$ws.Range("U8").Value = 339
$ws.Range("V8").Value = 75
$ws.Range("W8").Value = 8
$ws.Range("X8").Value = 9070
